$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 96
$ws.Range("A96").Value = 95.0
$ws.Range("B96").Value = 'Monday, Jan 09'
$ws.Range("C96").Value = '1:40 PM'
$ws.Range("D96").Value = 'LH1356'
$ws.Range("E96").Value = 'Frankfurt'
$ws.Range("F96").Value = '(FRA)'
$ws.Range("G96").Value = 'Lufthansa '
$ws.Range("H96").Value = 'CRJ9'
$ws.Range("I96").Value = '(D-ACNJ)'
$ws.Range("J96").Value = '1:54 PM'
$ws.Range("L96").Value = '0 hours, 14 minutes'

# Row 97
$ws.Range("A97").Value = 96.0
$ws.Range("B97").Value = 'Monday, Jan 09'
$ws.Range("C97").Value = '2:12 PM'
$ws.Range("D97").Value = 'UNKNOWN'
$ws.Range("E97").Value = 'Basel'
$ws.Range("F97").Value = '(BSL)'
$ws.Range("G97").Value = 'Enter Air '
$ws.Range("H97").Value = 'B738'
$ws.Range("I97").Value = '(SP-ENR)'
$ws.Range("J97").Value = '2:18 PM'
$ws.Range("L97").Value = '0 hours, 6 minutes'

# Row 98
$ws.Range("A98").Value = 97.0
$ws.Range("B98").Value = 'Monday, Jan 09'
$ws.Range("C98").Value = '2:40 PM'
$ws.Range("D98").Value = 'LO3883'
$ws.Range("E98").Value = 'Warsaw'
$ws.Range("F98").Value = '(WAW)'
$ws.Range("G98").Value = 'LOT (Sliwka Naleczowska Livery) '
$ws.Range("H98").Value = 'E195'
$ws.Range("I98").Value = '(SP-LNC)'
$ws.Range("J98").Value = '2:48 PM'
$ws.Range("L98").Value = '0 hours, 8 minutes'

# Row 99
$ws.Range("A99").Value = 98.0
$ws.Range("B99").Value = 'Monday, Jan 09'
$ws.Range("C99").Value = '4:10 PM'
$ws.Range("D99").Value = 'KL1815'
$ws.Range("E99").Value = 'Amsterdam'
$ws.Range("F99").Value = '(AMS)'
$ws.Range("G99").Value = 'KLM '
$ws.Range("H99").Value = 'E295'
$ws.Range("I99").Value = '(PH-NXF)'
$ws.Range("J99").Value = '3:59 PM'
$ws.Range("L99").Value = '0 hours, -11 minutes'

# Row 100
$ws.Range("A100").Value = 99.0
$ws.Range("B100").Value = 'Monday, Jan 09'
$ws.Range("C100").Value = '5:00 PM'
$ws.Range("D100").Value = 'FR7101'
$ws.Range("E100").Value = 'Oslo'
$ws.Range("F100").Value = '(OSL)'
$ws.Range("G100").Value = 'Ryanair '
$ws.Range("H100").Value = 'B738'
$ws.Range("I100").Value = '(SP-RSO)'
$ws.Range("J100").Value = '4:49 PM'
$ws.Range("L100").Value = '0 hours, -11 minutes'

# Row 101
$ws.Range("A101").Value = 100.0
$ws.Range("B101").Value = 'Monday, Jan 09'
$ws.Range("C101").Value = '5:25 PM'
$ws.Range("D101").Value = 'W61072'
$ws.Range("E101").Value = 'Eindhoven'
$ws.Range("F101").Value = '(EIN)'
$ws.Range("G101").Value = 'Wizz Air '
$ws.Range("H101").Value = 'A321'
$ws.Range("I101").Value = '(HA-LXP)'
$ws.Range("J101").Value = '5:23 PM'
$ws.Range("L101").Value = '0 hours, -2 minutes'

# Row 102
$ws.Range("A102").Value = 101.0
$ws.Range("B102").Value = 'Monday, Jan 09'
$ws.Range("C102").Value = '5:50 PM'
$ws.Range("D102").Value = 'W61012'
$ws.Range("E102").Value = 'Liverpool'
$ws.Range("F102").Value = '(LPL)'
$ws.Range("G102").Value = 'Wizz Air '
$ws.Range("H102").Value = 'A321'
$ws.Range("I102").Value = '(HA-LXD)'
$ws.Range("J102").Value = '5:28 PM'
$ws.Range("L102").Value = '0 hours, -22 minutes'
